# Update cryptocurrency price/volume data per upstream diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.081.02'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').Value = '3.124.86'
$ws.Range('E3').Value = '  +0.55%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '173.75'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.30%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.155'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.45%  '
$ws.Range('E10').Value = '  -1.82%  '
$ws.Range('E11').Value = '  -0.51%  '
$ws.Range('E12').Value = '  +0.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '37.32'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.20%  '
$ws.Range('E14').Value = '  -1.65%  '
$ws.Range('D15').Value = '3.644.79'
$ws.Range('E15').Value = '  +0.60%  '
$ws.Range('D16').Value = '67.068.28'
$ws.Range('E16').Value = '  +0.18%  '
$ws.Range('E17').Value = '  -0.71%  '
$ws.Range('D18').Value = '3.128.12'
$ws.Range('E18').Value = '  +0.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.41'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.76%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '492.91'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.54%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.93'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.41%  '
$ws.Range('B22').Value = 'Polygon'
$ws.Range('C22').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.708'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.01%  '
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.29'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.29'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.44'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.91'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.85%  '
$ws.Range('E29').Value = '  -2.11%  '
$ws.Range('E30').Value = '  -0.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '28.66'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.42%  '
$ws.Range('E32').Value = '  -0.67%  '
$ws.Range('D33').Value = '0.0₃0948'
$ws.Range('E33').Value = '  -6.62%  '
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.90'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.978'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '46.97'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.79%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.06'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.13%  '
$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.312'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.00%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.124'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.81%  '
$ws.Range('B41').Value = 'Cosmos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.56'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.53%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '2.839.17'
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.62'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.00%  '
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '385.02'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0353'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.43%  '
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '135.72'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.57%  '
$ws.Range('B47').Value = 'USDe'
$ws.Range('C47').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.00'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '25.01'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.03%  '
$ws.Range('B49').Value = 'ThetaToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.22'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.87%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.108'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.83%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.76'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.97%  '
